$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix self-referencing SUM ranges in the original demonstration block (rows 2-10) ---
# (these previously included the formula cell itself in its own SUM range)
$ws.Range("Q2").Formula = "=MOD(SUM(R2:T10,Q3:Q10),2)"
$ws.Range("M3").Formula = "=MOD(SUM(N3:T3,M5:T5,M7:T7,M9:T9),2)"
$ws.Range("M4").Formula = "=MOD(SUM(N4:T4,M5:T5,M8:T9),2)"
$ws.Range("M6").Formula = "=MOD(SUM(M7:T9,N6:T6),2)"

# --- Add the new worked example block (rows 12-21), mirroring the layout of rows 1-10 ---
# Row 12
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 4
$ws.Range("G12").Value = 5
$ws.Range("H12").Value = 6
$ws.Range("I12").Value = 7
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 2
$ws.Range("P12").Value = 3
$ws.Range("Q12").Value = 4
$ws.Range("R12").Value = 5
$ws.Range("S12").Value = 6
$ws.Range("T12").Value = 7

# Row 13
$ws.Range("A13").Value = 0
$ws.Range("C13").Formula = "=MOD(MOD(SUM(C14:C21,E13:E21,G13:G21,I13:I21),2),2)"
$ws.Range("D13").Formula = "=MOD(MOD(SUM(D14:E21,H13:I21,E13),2),2)"
$ws.Range("E13").Value = 1
$ws.Range("F13").Formula = "=MOD(SUM(G13:I21,F14:F21),2)"
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").Formula = "=MOD(MOD(SUM(N14:N21,P13:P21,R13:R21,T13:T21),2),2)"
$ws.Range("O13").Formula = "=MOD(MOD(SUM(O14:P21,S13:T21,P13),2),2)"
$ws.Range("P13").Value = 1
$ws.Range("Q13").Formula = "=MOD(SUM(R13:T21,Q14:Q21),2)"
$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 0
$ws.Range("T13").Value = 0

# Row 14
$ws.Range("A14").Value = 1
$ws.Range("B14").Formula = "=MOD(SUM(C14:I14,B16:I16,B18:I18,B20:I20),2)"
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("L14").Value = 1
$ws.Range("M14").Formula = "=MOD(SUM(N14:T14,M16:T16,M18:T18,M20:T20),2)"
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 1
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = 0
$ws.Range("S14").Value = 0
$ws.Range("T14").Value = 0

# Row 15
$ws.Range("A15").Value = 2
$ws.Range("B15").Formula = "=MOD(SUM(C15:I15,B16:I16,B19:I20),2)"
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 1
$ws.Range("L15").Value = 2
$ws.Range("M15").Formula = "=MOD(SUM(N15:T15,M16:T16,M19:T20),2)"
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 1
$ws.Range("P15").Value = 1
$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = 1
$ws.Range("S15").Value = 0
$ws.Range("T15").Value = 1

# Row 16
$ws.Range("A16").Value = 3
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 0
$ws.Range("L16").Value = 3
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 1
$ws.Range("P16").Value = 1
$ws.Range("Q16").Value = 0
$ws.Range("R16").Value = 0
$ws.Range("S16").Value = 1
$ws.Range("T16").Value = 0

# Row 17
$ws.Range("A17").Value = 4
$ws.Range("B17").Formula = "=MOD(SUM(B18:I20,C17:I17),2)"
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = 0
$ws.Range("L17").Value = 4
$ws.Range("M17").Formula = "=MOD(SUM(M18:T20,N17:T17),2)"
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 1
$ws.Range("Q17").Value = 0
$ws.Range("R17").Value = 1
$ws.Range("S17").Value = 1
$ws.Range("T17").Value = 0

# Row 18
$ws.Range("A18").Value = 5
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 1
$ws.Range("I18").Value = 1
$ws.Range("L18").Value = 5
$ws.Range("M18").Value = 1
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 1
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = 0
$ws.Range("S18").Value = 1
$ws.Range("T18").Value = 1

# Row 19
$ws.Range("A19").Value = 6
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 1
$ws.Range("I19").Value = 0
$ws.Range("L19").Value = 6
$ws.Range("M19").Value = 1
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 1
$ws.Range("R19").Value = 0
$ws.Range("S19").Value = 1
$ws.Range("T19").Value = 0

# Row 20
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 0
$ws.Range("L20").Value = 7
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 1
$ws.Range("Q20").Value = 0
$ws.Range("R20").Value = 1
$ws.Range("S20").Value = 1
$ws.Range("T20").Value = 0

# Row 21
$ws.Range("A21").Value = 8
$ws.Range("B21").Formula = "=MOD(SUM(B21:I21),2)"
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("L21").Value = 8
$ws.Range("M21").Formula = "=MOD(SUM(M21:T21),2)"
$ws.Range("N21").Value = 0
$ws.Range("O21").Value = 1
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 1
$ws.Range("R21").Value = 1
$ws.Range("S21").Value = 0
$ws.Range("T21").Value = 0

# --- Apply the red-font / yellow-fill "formula" style used elsewhere in the sheet ---
# (matches the style already used on N2, O2, Q2, M3, M4, M6, M10, etc.)
$ws.Range("N2").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("F13").PasteSpecial(-4122)
$ws.Range("N13").PasteSpecial(-4122)
$ws.Range("O13").PasteSpecial(-4122)
$ws.Range("Q13").PasteSpecial(-4122)
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("M14").PasteSpecial(-4122)
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("M15").PasteSpecial(-4122)
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("M17").PasteSpecial(-4122)
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("M21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the active selection to match the author's final cursor position ---
$ws.Range("P17").Select()